# Add a new row (19) to the Assignments tracker sheet for "Assignment_18",
# matching the style/formatting of the preceding row (18), plus its hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content ---------------------------------------------------
$ws.Range("A19").Value = "Assignment_18"
$ws.Range("B19").Value = "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_18/CRUD_API"
$ws.Range("C19").Value = 45174

# --- Hyperlink for the new URL cell ------------------------------------
$ws.Hyperlinks.Add($ws.Range("B19"), "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_18/CRUD_API")

# --- Match formatting/styles of the row above (row 18) ------------------
$ws.Range("A18:C18").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
[void]$ws.Application.CutCopyMode

# --- Row height to match the other data rows -----------------------------
$ws.Rows.Item(19).RowHeight = 32.25

# --- Update the active selection to the cell after the new row (D19) -----
[void]$ws.Range("D19").Select()

Write-Host "Added Assignment_18 row"
